# "Work on shop (3) & improve map"
# - Adds PRICE / GOLD DROP columns (U, V) to the CHAMPION sheet with values per champion.
# - Tweaks the upgrade-step count (column T) for two rows.
# - Re-does the frozen-pane / selection layout on CHAMPION (freeze col A + row 1).
# - Updates the ITEM sheet's scroll position / selection and makes CHAMPION the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHAMPION")
$itemWs = $wb.Worksheets.Item("ITEM")

# --- New headers: PRICE (U1) and GOLD DROP (V1) ---
$ws.Range("U1").Value = "PRICE"
$ws.Range("V1").Value = "GOLD DROP"
# Match header style/formatting used by the rest of row 1
$ws.Range("U1").Style = $ws.Range("T1").Style
$ws.Range("V1").Style = $ws.Range("T1").Style

# Give the new GOLD DROP column (V / col 22) a sensible width, same as the source edit
$ws.Columns.Item(22).ColumnWidth = 11.6640625

# Header row now wraps onto two lines given the narrower new column
$ws.Rows.Item(1).RowHeight = 28.8

# --- PRICE / GOLD DROP values per champion row ---
$priceGold = @{
    3  = @(60, 15)
    4  = @(65, 13)
    5  = @(200, 50)
    7  = @(125, 35)
    9  = @(70, 20)
    10 = @(55, 15)
    11 = @(50, 19)
    12 = @(100, 33)
    14 = @(58, 14)
    15 = @(120, 37)
    16 = @(60, 18)
    17 = @(62, 16)
    18 = @(52, 16)
    19 = @(115, 36)
    20 = @(54, 15)
    21 = @(215, 50)
    24 = @(64, 15)
    25 = @(67, 14)
    26 = @(350, 95)
}

foreach ($row in $priceGold.Keys) {
    $vals = $priceGold[$row]
    $ws.Range("U$row").Value = $vals[0]
    $ws.Range("V$row").Value = $vals[1]
    $ws.Range("U$row").Style = $ws.Range("T$row").Style
    $ws.Range("V$row").Style = $ws.Range("T$row").Style
}

# --- Upgrade-step count tweaks ---
$ws.Range("T4").Value = 1
$ws.Range("T26").Value = 4

# --- CHAMPION sheet view: freeze header row + first column, select W24 ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("W24").Select()

# --- ITEM sheet view: scroll position + selection, no longer the active tab ---
$itemWs.Activate()
$excel.ActiveWindow.ScrollRow = 8
$itemWs.Range("E2").Select()

# Leave CHAMPION as the active/selected sheet, matching the source edit
$ws.Activate()
